$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gof_FlxExt")

# Fill in the previously empty coefficient columns (E:J) for rows 7-10
$ws.Range("E7").Value = 4.5197081384539599
$ws.Range("F7").Value = 3.9722004305122902
$ws.Range("G7").Value = 11.7412109674933
$ws.Range("H7").Value = 2.0417894891824799
$ws.Range("I7").Value = 0.84442496918497401
$ws.Range("J7").Value = 3.58761452437188

$ws.Range("E8").Value = 1.99920194837266
$ws.Range("F8").Value = 0.89213700366721804
$ws.Range("G8").Value = 4.6092079774300299
$ws.Range("H8").Value = 0.89562322843932496
$ws.Range("I8").Value = 0.17904777063236799
$ws.Range("J8").Value = 1.8463910818920899

$ws.Range("E9").Value = 0.84470560285347895
$ws.Range("F9").Value = 0.23173900510265799
$ws.Range("G9").Value = 1.95145372442261
$ws.Range("H9").Value = 0.68178489628790195
$ws.Range("I9").Value = 0.150967387789747
$ws.Range("J9").Value = 2.0676835787835999

$ws.Range("E10").Value = 3.03222221811638
$ws.Range("F10").Value = 4.4855531007687901
$ws.Range("G10").Value = 5.5870357213666804
$ws.Range("H10").Value = 1.4609329694573101
$ws.Range("I10").Value = 1.0412488414276799
$ws.Range("J10").Value = 2.7551515923112402

# Update the selected cell on this sheet to match the saved selection
$ws.Activate()
$ws.Range("I12").Select()
